# Extend the "group" answer columns by one more column (BY), mirroring the
# existing last column (BX) for every data row, as in the commit that added
# an extra duplicated answer column to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$src = $ws.Range("BX2:BX15")
$dst = $ws.Range("BY2:BY15")

$src.Copy($dst)
